$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update F column "想去人数" values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2133
$ws1.Range("F4").Value = 20
$ws1.Range("F5").Value = 11103
$ws1.Range("F10").Value = 10998
$ws1.Range("F14").Value = 1713
$ws1.Range("F15").Value = 5508
$ws1.Range("F17").Value = 3425

# Sheet "全部类型" (sheet4) - update F column "想去人数" values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2133
$ws4.Range("F5").Value = 20
$ws4.Range("F7").Value = 11103
$ws4.Range("F12").Value = 10998
$ws4.Range("F16").Value = 1713
$ws4.Range("F17").Value = 5508
$ws4.Range("F19").Value = 3425
